# Meeting-auto-attender workbook refresh: swap in the next meeting's
# date/time + Zoom link, and blank out the now-stale Meeting ID / password.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- A2: new meeting date/time, entered as text (leading apostrophe) ---
# Re-apply the existing date display format afterwards so the cell keeps
# its original "short date" number format instead of picking up a fresh
# custom one.
$dateCell = $ws.Range("A2")
$dateCell.Formula = "'20-08-2021 12:00 PM"
$dateCell.NumberFormat = "mm-dd-yy"

# --- B2: new Zoom link text (hyperlink target itself is left as-is) ---
$ws.Range("B2").Value = "https://us02web.zoom.us/j/85071211231231231"

# --- C2: clear the stale Meeting ID, keep its formatting ---
$ws.Range("C2").ClearContents()

# --- D2: clear the stale Meeting password entirely (content + format) ---
$ws.Range("D2").Clear()

# Leave the selection where the author last left it before saving.
$ws.Range("C26").Select() | Out-Null
